# Generate Report for Handoff
# Updates the localization-status workbook: the 9cc7a019 and f0495d95 source
# files have been re-handed-off, so their status flips from
# "Handed back: in sync with en-US" to "Ready for handoff", their
# handoff timestamps move forward, and an "Error Detail" note is recorded
# (handback file version is stale) on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$overviewDate = "2016-08-31 02:29:22"
$zhcnHandoffDate = "2016-08-31 02:29:17"
$dedeHandoffDate = "2016-08-31 02:29:22"

$err9cc7 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a0f807ec86aafdb55e563267bb96843a7481e339/e2e/9cc7a019-9efc-47a8-be64-c4dc30c3600d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01c1dc7359e15549b90a2aed78925d881d688ee8/e2e/9cc7a019-9efc-47a8-be64-c4dc30c3600d.md."
$errf049 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a0f807ec86aafdb55e563267bb96843a7481e339/e2e/f0495d95-1400-47db-a186-4b4276ea9238.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01c1dc7359e15549b90a2aed78925d881d688ee8/e2e/f0495d95-1400-47db-a186-4b4276ea9238.md."

# --- Overview sheet: per-locale status columns (E=zh-cn, F=de-de) + latest
#     handoff-xliff-generate date (G), for rows 4 (9cc7a019) and 5 (f0495d95)
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = $statusReady
$wsOverview.Range("F4").Value = $statusReady
$wsOverview.Range("G4").Value = $overviewDate

$wsOverview.Range("E5").Value = $statusReady
$wsOverview.Range("F5").Value = $statusReady
$wsOverview.Range("G5").Value = $overviewDate

# --- zh-cn sheet: Status (C), Latest Handoff Datetime (H), Error Detail (P)
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $statusReady
$wsZhCn.Range("H4").Value = $zhcnHandoffDate
$wsZhCn.Range("P4").Value = $err9cc7

$wsZhCn.Range("C5").Value = $statusReady
$wsZhCn.Range("H5").Value = $zhcnHandoffDate
$wsZhCn.Range("P5").Value = $errf049

# Error Detail column got wider to fit the new message
$wsZhCn.Range("P1").EntireColumn.ColumnWidth = 39.166666666666664

# --- de-de sheet: Status (C), Latest Handoff Datetime (H), Error Detail (P)
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $statusReady
$wsDeDe.Range("H4").Value = $dedeHandoffDate
$wsDeDe.Range("P4").Value = $err9cc7

$wsDeDe.Range("C5").Value = $statusReady
$wsDeDe.Range("H5").Value = $dedeHandoffDate
$wsDeDe.Range("P5").Value = $errf049

$wsDeDe.Range("P1").EntireColumn.ColumnWidth = 39.166666666666664
